# Update "想去人数" (interest count) values on the "展览" and "全部类型" sheets
# to reflect the latest scrape output, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 315
$wsExhibit.Range("F5").Value = 2825
$wsExhibit.Range("F6").Value = 1968
$wsExhibit.Range("F7").Value = 378
$wsExhibit.Range("F8").Value = 129
$wsExhibit.Range("F9").Value = 1026
$wsExhibit.Range("F10").Value = 192
$wsExhibit.Range("F11").Value = 150
$wsExhibit.Range("F12").Value = 35

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 315
$wsAll.Range("F5").Value = 2825
$wsAll.Range("F6").Value = 1968
$wsAll.Range("F7").Value = 378
$wsAll.Range("F9").Value = 129
$wsAll.Range("F10").Value = 1026
$wsAll.Range("F11").Value = 192
$wsAll.Range("F12").Value = 150
$wsAll.Range("F13").Value = 35
